# Scheduled-runner style update of market-board derived profit columns (H:N)
# across the per-job Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only numeric value cells are touched; no formulas/styles in this workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3720
$ws.Range("I64").Value = 5266.6665
$ws.Range("J64").Value = 3057.1428
$ws.Range("K64").Value = 5266.6665
$ws.Range("L64").Value = 3057.1428
$ws.Range("M64").Value = -5018.6665
$ws.Range("N64").Value = -3553.1428

$ws.Range("H67").Value = 3720
$ws.Range("I67").Value = 5266.6665
$ws.Range("J67").Value = 3057.1428
$ws.Range("K67").Value = 5266.6665
$ws.Range("L67").Value = 3057.1428
$ws.Range("M67").Value = -4408.6665
$ws.Range("N67").Value = -4773.1428

$ws.Range("H103").Value = 1019.0909
$ws.Range("I103").Value = 1041
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 3123
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -2537
$ws.Range("N103").Value = -3572

$ws.Range("H131").Value = 1788.1666
$ws.Range("I131").Value = 1062.2858
$ws.Range("J131").Value = 4328.75
$ws.Range("K131").Value = 3186.8574
$ws.Range("L131").Value = 12986.25
$ws.Range("M131").Value = 1853.1426
$ws.Range("N131").Value = -23066.25

$ws.Range("H135").Value = 732.6515000000001
$ws.Range("I135").Value = 403.75
$ws.Range("J135").Value = 1954.2858
$ws.Range("K135").Value = 3633.75
$ws.Range("L135").Value = 17588.5722
$ws.Range("M135").Value = -1098.75
$ws.Range("N135").Value = -22658.5722

$ws.Range("H137").Value = 893
$ws.Range("I137").Value = 850.3333
$ws.Range("J137").Value = 907.2222
$ws.Range("K137").Value = 2550.9999
$ws.Range("L137").Value = 2721.6666
$ws.Range("M137").Value = -0.9998999999997977
$ws.Range("N137").Value = -7821.6666

$ws.Range("H138").Value = 1572.35
$ws.Range("I138").Value = 759.7193
$ws.Range("J138").Value = 2649.558
$ws.Range("K138").Value = 2279.1579
$ws.Range("L138").Value = 7948.674
$ws.Range("M138").Value = 2860.8421
$ws.Range("N138").Value = -18228.674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 16867.75
$ws.Range("I28").Value = 9823.666999999999
$ws.Range("J28").Value = 38000
$ws.Range("K28").Value = 9823.666999999999
$ws.Range("L28").Value = 38000
$ws.Range("M28").Value = -9631.666999999999
$ws.Range("N28").Value = -38384

$ws.Range("H32").Value = 1430.64
$ws.Range("I32").Value = 1210.2317
$ws.Range("J32").Value = 2434.7222
$ws.Range("K32").Value = 1210.2317
$ws.Range("L32").Value = 2434.7222
$ws.Range("M32").Value = -923.2317
$ws.Range("N32").Value = -3008.7222

$ws.Range("H61").Value = 976.0909
$ws.Range("J61").Value = 1478.2222
$ws.Range("L61").Value = 1478.2222
$ws.Range("N61").Value = -1902.2222

$ws.Range("H74").Value = 1264.5
$ws.Range("I74").Value = 1231
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 1231
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = -357
$ws.Range("N74").Value = -3448

$ws.Range("H77").Value = 1264.5
$ws.Range("I77").Value = 1231
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 6155
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = -1787
$ws.Range("N77").Value = -17236

$ws.Range("H97").Value = 537.8946999999999
$ws.Range("I97").Value = 501.1111
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 501.1111
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -5.111100000000022
$ws.Range("N97").Value = -2192

$ws.Range("H99").Value = 16867.75
$ws.Range("I99").Value = 9823.666999999999
$ws.Range("J99").Value = 38000
$ws.Range("K99").Value = 9823.666999999999
$ws.Range("L99").Value = 38000
$ws.Range("M99").Value = -6828.666999999999
$ws.Range("N99").Value = -43990

$ws.Range("H102").Value = 4625
$ws.Range("I102").Value = 4833.3335
$ws.Range("K102").Value = 4833.3335
$ws.Range("M102").Value = -3211.3335

$ws.Range("H136").Value = 976.0909
$ws.Range("J136").Value = 1478.2222
$ws.Range("L136").Value = 4434.6666
$ws.Range("N136").Value = -9534.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 62502090
$ws.Range("I99").Value = 83335420
$ws.Range("J99").Value = 2080.5
$ws.Range("K99").Value = 83335420
$ws.Range("L99").Value = 2080.5
$ws.Range("M99").Value = -83333922
$ws.Range("N99").Value = -5076.5

$ws.Range("H134").Value = 17785
$ws.Range("I134").Value = 1388.8868
$ws.Range("J134").Value = 114339.89
$ws.Range("K134").Value = 4166.6604
$ws.Range("L134").Value = 343019.67
$ws.Range("M134").Value = -1631.6604
$ws.Range("N134").Value = -348089.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2377.9016
$ws.Range("I31").Value = 2231
$ws.Range("J31").Value = 2638.318
$ws.Range("K31").Value = 2231
$ws.Range("L31").Value = 2638.318
$ws.Range("M31").Value = -1936
$ws.Range("N31").Value = -3228.318

$ws.Range("H34").Value = 2377.9016
$ws.Range("I34").Value = 2231
$ws.Range("J34").Value = 2638.318
$ws.Range("K34").Value = 2231
$ws.Range("L34").Value = 2638.318
$ws.Range("M34").Value = -2029
$ws.Range("N34").Value = -3042.318

$ws.Range("H58").Value = 2946.92
$ws.Range("I58").Value = 1041.8064
$ws.Range("J58").Value = 6055.263
$ws.Range("K58").Value = 1041.8064
$ws.Range("L58").Value = 6055.263
$ws.Range("M58").Value = -838.8063999999999
$ws.Range("N58").Value = -6461.263

$ws.Range("H107").Value = 359.75
$ws.Range("I107").Value = 356.8
$ws.Range("J107").Value = 361.85715
$ws.Range("K107").Value = 356.8
$ws.Range("L107").Value = 361.85715
$ws.Range("M107").Value = 1563.2
$ws.Range("N107").Value = -4201.85715

$ws.Range("H132").Value = 1502.5508
$ws.Range("I132").Value = 913.0952
$ws.Range("K132").Value = 2739.2856
$ws.Range("M132").Value = -209.2856000000002

$ws.Range("H136").Value = 2946.92
$ws.Range("I136").Value = 1041.8064
$ws.Range("J136").Value = 6055.263
$ws.Range("K136").Value = 3125.4192
$ws.Range("L136").Value = 18165.789
$ws.Range("M136").Value = -575.4191999999998
$ws.Range("N136").Value = -23265.789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 6014.5
$ws.Range("J104").Value = 6014.5
$ws.Range("L104").Value = 18043.5
$ws.Range("N104").Value = -23285.5

$ws.Range("H107").Value = 648817.3
$ws.Range("I107").Value = 729
$ws.Range("K107").Value = 2187
$ws.Range("M107").Value = -267

$ws.Range("H129").Value = 4639.125
$ws.Range("I129").Value = 1690
$ws.Range("J129").Value = 5622.1665
$ws.Range("K129").Value = 5070
$ws.Range("L129").Value = 16866.4995
$ws.Range("M129").Value = -70
$ws.Range("N129").Value = -26866.4995

$ws.Range("H137").Value = 39335.25
$ws.Range("I137").Value = 1953.5294
$ws.Range("J137").Value = 97107
$ws.Range("K137").Value = 5860.5882
$ws.Range("L137").Value = 291321
$ws.Range("M137").Value = -760.5882000000001
$ws.Range("N137").Value = -301521

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3966.5454
$ws.Range("I70").Value = 3829.7437
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 3829.7437
$ws.Range("L70").Value = 4300
$ws.Range("M70").Value = -3559.7437
$ws.Range("N70").Value = -4840

$ws.Range("H73").Value = 3966.5454
$ws.Range("I73").Value = 3829.7437
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 3829.7437
$ws.Range("L73").Value = 4300
$ws.Range("M73").Value = -2893.7437
$ws.Range("N73").Value = -6172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2826
$ws.Range("I68").Value = 3530
$ws.Range("J68").Value = 1770
$ws.Range("K68").Value = 3530
$ws.Range("L68").Value = 1770
$ws.Range("M68").Value = -2781
$ws.Range("N68").Value = -3268

$ws.Range("H71").Value = 2826
$ws.Range("I71").Value = 3530
$ws.Range("J71").Value = 1770
$ws.Range("K71").Value = 17650
$ws.Range("L71").Value = 8850
$ws.Range("M71").Value = -13906
$ws.Range("N71").Value = -16338

$ws.Range("H136").Value = 2297.465
$ws.Range("I136").Value = 1344.6552
$ws.Range("K136").Value = 4033.9656
$ws.Range("M136").Value = -1483.9656

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5016.6665
$ws.Range("I100").Value = 6250
$ws.Range("J100").Value = 4400
$ws.Range("K100").Value = 12500
$ws.Range("L100").Value = 8800
$ws.Range("M100").Value = -11959
$ws.Range("N100").Value = -9882

$ws.Range("H132").Value = 473.4304
$ws.Range("I132").Value = 458.27692
$ws.Range("J132").Value = 543.7857
$ws.Range("K132").Value = 1374.83076
$ws.Range("L132").Value = 1631.3571
$ws.Range("M132").Value = 1155.16924
$ws.Range("N132").Value = -6691.3571

$ws.Range("H136").Value = 741.76
$ws.Range("I136").Value = 898.0333000000001
$ws.Range("J136").Value = 507.35
$ws.Range("K136").Value = 2694.0999
$ws.Range("L136").Value = 1522.05
$ws.Range("M136").Value = -144.0999000000002
$ws.Range("N136").Value = -6622.05
